$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 8029.442592859268
$ws.Range("B2").Value = 133.8240432143211
$ws.Range("C2").Value = 729.586681951176
$ws.Range("D2").Value = 11
$ws.Range("E2").Value = 0.08219748660846647
